$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (row-by-row, matches upstream scrape refresh)

# Row 2
$ws.Range('D2').Value = '90.929.11'
$ws.Range('E2').Value = '  -4.31%  '

# Row 3
$ws.Range('D3').Value = '3.295.00'
$ws.Range('E3').Value = '  -5.53%  '

# Row 4
$ws.Range('E4').Value = '  -0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '228.57'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -4.50%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '610.63'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -5.25%  '

# Row 7
$ws.Range('E7').Value = '  -5.34%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.380'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -6.24%  '

# Row 9
$ws.Range('E9').Value = '  +0.04%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.940'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -6.25%  '

# Row 11
$ws.Range('D11').Value = '3.289.98'
$ws.Range('E11').Value = '  -5.74%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '41.51'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -3.20%  '

# Row 13
$ws.Range('E13').Value = '  -3.53%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.92'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -4.75%  '

# Row 15
$ws.Range('D15').Value = '90.821.47'
$ws.Range('E15').Value = '  -4.31%  '

# Row 16
$ws.Range('D16').Value = '3.915.11'
$ws.Range('E16').Value = '  -5.67%  '

# Row 17
$ws.Range('E17').Value = '  -6.27%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '8.00'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -5.99%  '

# Row 19
$ws.Range('D19').Value = '3.298.40'
$ws.Range('E19').Value = '  -5.49%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.00'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -5.77%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.77'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -6.02%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.37'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +5.29%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '484.33'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -4.47%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.438'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -14.20%  '

# Row 25
$ws.Range('E25').Value = '  -7.43%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.04'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -9.97%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '88.78'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -6.95%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '11.67'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -4.34%  '

# Row 29
$ws.Range('D29').Value = '3.479.79'
$ws.Range('E29').Value = '  -5.41%  '

# Row 30
$ws.Range('E30').Value = '  +0.03%  '

# Row 31
$ws.Range('E31').Value = '  -8.47%  '

# Row 32
$ws.Range('E32').Value = '  -1.61%  '

# Row 33
$ws.Range('E33').Value = '  -5.99%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.36%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.169'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -8.24%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '27.79'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -10.25%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.517'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -9.91%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '539.88'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -3.11%  '

# Row 39
$ws.Range('E39').Value = '  -0.02%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '7.25'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -6.91%  '

# Row 41
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.145'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.82%  '

# Row 42
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.35'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -8.39%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.852'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -9.61%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '23.66'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.75%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.62'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +2.23%  '

# Row 46
$ws.Range('B46').Value = 'ImmutableX'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.64'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.50%  '

# Row 47
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0404'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.89%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '5.31'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -6.88%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.07'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -4.87%  '

# Row 50
$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '51.36'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.94%  '

# Row 51
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.85'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -3.11%  '
